$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Author line: "Bartomeus, I., Stavert, J.R., Ward, D., and Aguado, O."
#    (unchanged visible text; upstream diff only adds proofErr spell-check
#    markers / run-splits around "Stavert" and "Aguado" -- no visible text
#    change, so nothing to do here.)
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2. Remove the stray "_GoBack" bookmark after "...different row."
# ---------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3. "Gbif" paragraph -- no visible text change (just proofErr wrapping).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4. "All files and data will be available at" -> "All files and data are available at"
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("All files and data will be available at", $true, $false, $false, $false, $false, $true, 1, $false, "All files and data are available at", 2)

# ---------------------------------------------------------------------
# 5. Replace the tail of the Sup Mat 3 paragraph:
#    " and deposited on Dryad or Figshare upon acceptance."
#    ->
#    " (DOI: 10.5281/zenodo.1326309)."
#    with a real hyperlink around the DOI, and re-insert the "_GoBack"
#    bookmark right after the DOI hyperlink.
# ---------------------------------------------------------------------
$tail = $d.Content.Find.Execute("and deposited on Dryad or Figshare upon acceptance.")
$hl = $d.Hyperlinks(1)
$insertPoint = $hl.Range.End
$r = $d.Range($insertPoint, $insertPoint)
# find and remove the old trailing sentence first
$oldTail = $d.Range($insertPoint, $d.Paragraphs(9).Range.End)
$oldTail.Text = ""

$cur = $insertPoint
$r1 = $d.Range($cur, $cur)
$r1.InsertAfter(" (DOI: ")
$cur = $cur + " (DOI: ".Length

$doiText = "10.5281/zenodo.1326309"
$r2 = $d.Range($cur, $cur)
$r2.InsertAfter($doiText)
$doiRange = $d.Range($cur, $cur + $doiText.Length)
$d.Hyperlinks.Add($doiRange, "https://doi.org/10.5281/zenodo.1326309", "", "", $doiText) | Out-Null
$cur = $cur + $doiText.Length

$d.Bookmarks.Add("_GoBack", $d.Range($cur, $cur)) | Out-Null

$r3 = $d.Range($cur, $cur)
$r3.InsertAfter(").")
